# Redid tables and plots
# Update the keyterm co-occurrence table (1994-2003): refresh cluster
# key terms / sizes / centrality / density, drop the old "emotion ..."
# cluster row, and promote the former cluster 6 row (eeg ...) into the
# now-empty cluster-5 slot with its own refreshed numbers/colour.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# wdColor value for RGB hex AABBCC is R + G*256 + B*65536
$purple = 12412820   # 9467BD

# --- Remove the old cluster-5 row ("emotion, affect, startle, ...") ---
$t.Rows.Item(6).Delete()

# --- Cluster 1 (row 2) ---
$t.Cell(2,2).Range.Text = "erp, human, memory, visual, auditory, electrophysiological, index, method, mmn, signal, temporal, detection"
$t.Cell(2,3).Range.Text = "12"
$t.Cell(2,4).Range.Text = "2261 (2)"
$t.Cell(2,5).Range.Text = "1699 (4)"

# --- Cluster 2 (row 3) ---
$t.Cell(3,2).Range.Text = "attention, eeg, child, development, sleep, adult, difference, individual, word, age, normal, pain"
$t.Cell(3,3).Range.Text = "12"
$t.Cell(3,4).Range.Text = "2297 (1)"
$t.Cell(3,5).Range.Text = "1937 (2)"

# --- Cluster 3 (row 4) ---
$t.Cell(4,2).Range.Text = "cardiovascular, heart rate, autonomic, blood pressure, stress, pattern, behavior, control, cardiac, respiratory"
$t.Cell(4,3).Range.Text = "10"
$t.Cell(4,4).Range.Text = "1668 (4)"
$t.Cell(4,5).Range.Text = "2068 (1)"

# --- Cluster 4 (row 5) ---
$t.Cell(5,2).Range.Text = "startle, emotion, brain, affect, perception, context, probe, complex, activation"
$t.Cell(5,3).Range.Text = "9"
$t.Cell(5,4).Range.Text = "2009 (3)"
$t.Cell(5,5).Range.Text = "1790 (3)"

# --- Former cluster 6 row (now row 6) becomes the new cluster 5 ---
$t.Cell(6,1).Range.Text = "5"
$t.Cell(6,2).Range.Text = "p300, scene, stimulus, patients, amplitude, anticipation, saccade"
$t.Cell(6,3).Range.Text = "7"
$t.Cell(6,4).Range.Text = "1394 (5)"
$t.Cell(6,5).Range.Text = "1640 (5)"

# Recolor the text runs only (exclude the trailing cell/paragraph mark so
# the paragraph-mark run properties are left untouched).
for ($col = 1; $col -le 5; $col++) {
    $cellRange = $t.Cell(6, $col).Range
    $textRange = $d.Range($cellRange.Start, $cellRange.End - 1)
    $textRange.Font.Color = $purple
}
